$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D28", "D29", "D30", "D34", "D35", "D36", "D38", "D40", "D41", "D44", "D46", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '60.946.54'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '3.386.23'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").Value = '571.62'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '142.34'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D8").Value = '0.474'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '7.62'
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  -1.87%  '
$ws.Range("D11").Value = '0.389'
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("D12").Value = '3.966.52'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '0.126'
$ws.Range("E13").Value = '  +1.81%  '
$ws.Range("D14").Value = '27.79'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000171'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.390.19'
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").Value = '61.069.84'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '6.09'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("D19").Value = '13.64'
$ws.Range("E19").Value = '  -4.39%  '
$ws.Range("D20").Value = '8.99'
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").Value = '382.57'
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").Value = '74.66'
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("D23").Value = '0.552'
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("E25").Value = '  -5.03%  '
$ws.Range("D26").Value = '3.526.74'
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").Value = '7.30'
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").Value = '8.01'
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("E32").Value = '  -4.42%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '23.30'
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("D35").Value = '6.96'
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("D36").Value = '167.28'
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.418.84'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '5.01'
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("E39").Value = '  -4.34%  '
$ws.Range("D40").Value = '0.0768'
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("D41").Value = '26.83'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("D44").Value = '4.38'
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").Value = '1.13'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '2.453.94'
$ws.Range("E47").Value = '  -4.79%  '
$ws.Range("D48").Value = '23.11'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").Value = '6.73'
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("D50").Value = '2.17'
$ws.Range("E50").Value = '  +8.39%  '
$ws.Range("D51").Value = '0.0264'
$ws.Range("E51").Value = '  +1.33%  '
